$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-28 Tuesday" "2025-10-29 Wednesday"

Replace-Text "28×62=1736" "80×80=6400"
Replace-Text "65×37=2405" "72×50=3600"
Replace-Text "34×54=1836" "32×80=2560"
Replace-Text "81×44=3564" "16×41=656"
Replace-Text "43×22=946" "63×39=2457"

Replace-Text "73×98=7154" "25×24=600"
Replace-Text "22×57=1254" "50×75=3750"
Replace-Text "39×81=3159" "46×79=3634"
Replace-Text "27×12=324" "79×56=4424"
Replace-Text "38×59=2242" "97×77=7469"

Replace-Text "14×49=686" "19×65=1235"
Replace-Text "37×40=1480" "28×20=560"
Replace-Text "66×82=5412" "63×81=5103"
Replace-Text "53×81=4293" "20×72=1440"
Replace-Text "74×51=3774" "84×51=4284"

Replace-Text "72×25=1800" "41×95=3895"
Replace-Text "25×72=1800" "38×61=2318"
Replace-Text "50×20=1000" "38×48=1824"
Replace-Text "74×69=5106" "74×82=6068"
Replace-Text "59×36=2124" "90×84=7560"

Replace-Text "39×14=546" "78×36=2808"
Replace-Text "32×30=960" "11×89=979"
Replace-Text "91×39=3549" "96×43=4128"
Replace-Text "16×97=1552" "77×69=5313"
Replace-Text "63×86=5418" "28×81=2268"
